{"js": "// The paragraph describing the supported file formats originally had its\n// bold run split in two (\"DOCX, DOC, PDF, HTML, XPS, R\" / \"TF and TXT\")\n// around a leftover \"_GoBack\" bookmark. Merge them back into a single run\n// (\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\") and drop the bookmark.\nconst body = context.document.body;\n\nconst results = body.search(\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replacing the found range (which spans the two runs + bookmark) with\n  // the same text merges it into a single run that keeps the original\n  // (bold) formatting of the first run.\n  results.items[0].insertText(\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\", \"Replace\");\n}\n\n// Remove the now-unneeded \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The paragraph describing the supported file formats originally had its\n# bold run split in two (\"DOCX, DOC, PDF, HTML, XPS, R\" / \"TF and TXT\")\n# around a leftover \"_GoBack\" bookmark. Merge them back into a single run\n# (\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\") and drop the bookmark.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Replacement.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n# wdFindContinue = 1, wdReplaceAll = 2 -- replacing the found range (which\n# spans the two runs + bookmark) with identical text merges it back into a\n# single run, keeping the original (bold) run formatting.\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Remove the now-unneeded \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
